# 2021-05-12 첨부파일 샘플 commit
# Insert a new "상세주소" (detailed address) column right before the
# existing "금융상품유형" column (column F), shifting the subsequent
# columns one place to the right, and populate the sample data for the
# new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; everything from F onward (values,
# styles, shared formulas, etc.) shifts right by one column.
$ws.Columns("F").Insert()

# The new column inherits the same width as the address column to its
# left (column E).
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# Header for the new column.
$ws.Range("F1").Value = "상세주소"

# Sample detail-address value for the data row. The guidance row (row 2)
# is intentionally left blank in the new column.
$ws.Range("F3").Value = "용산구 한남동 00-00"

# Reflect where the author's cursor ended up after the edit.
$ws.Range("F4").Select()
